# Weekly data refresh: insert a new price-report row for "Perejil" (Primera)
# at the top of the dated series (row 24), pushing the existing rows down by
# one. Excel's row-insert shifts all the data below it automatically, so the
# remaining rows (old 24..44 -> new 25..45) keep their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 24; everything currently at/after row 24
# (through the last used row, 44) moves down to make room.
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with this week's record.
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 45033
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = 100112044
$ws.Cells.Item(24, 7).Value = "Perejil"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 300
$ws.Cells.Item(24, 11).Value = 1500
$ws.Cells.Item(24, 12).Value = 1500
$ws.Cells.Item(24, 13).Value = 1500
$ws.Cells.Item(24, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(24, 15).Value = "Región del Maule"
$ws.Cells.Item(24, 16).Value = 1500
$ws.Cells.Item(24, 17).Value = 1
$ws.Cells.Item(24, 18).Value = "Hortaliza"
# Note: Rows.Item(24).Insert() already carries down the "Fecha" column's
# date number-format style (from the row that used to be at 24) onto the
# new row, so no explicit style copy is needed here.
